$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.816299999999999
$ws.Range("B10").Value = 5.092500000000002
$ws.Range("B12").Value = 5.000500000000001
$ws.Range("B18").Value = 6.980399999999995
$ws.Range("B37").Value = 8.631500000000003
$ws.Range("B55").Value = 6.308899999999995
$ws.Range("B68").Value = 4.874799999999996
$ws.Range("B77").Value = 9.009400000000007
$ws.Range("B78").Value = 9.360000000000003
$ws.Range("B81").Value = 5.419000000000003
$ws.Range("B82").Value = 5.313200000000002
